$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: D7/E7/F7 were stored as text ("45","1200","2500"); convert to real numbers ---
$ws.Range("D7").Value = 45
$ws.Range("E7").Value = 1200
$ws.Range("F7").Value = 2500
$ws.Range("G7").Value = 45808.680728125

# --- Row 10 (new): add the Dolex Gripa record ---
$ws.Range("A10").Value = "1M"
$ws.Range("B10").Value = "Medicamentos"
$ws.Range("C10").Value = "Dolex Gripa x 12 pastillas"

# D10/E10/F10 keep the quirky "numeric text" formatting used by the rest of
# the sheet's freshly-imported rows (inline/shared string holding a numeric
# look-alike, default style).
$ws.Range("D10:F10").NumberFormat = "@"
$ws.Range("D10").Value = "50"
$ws.Range("E10").Value = "4500"
$ws.Range("F10").Value = "6600"
$ws.Range("D10:F10").Style = $ws.Range("D9:F9").Style

# G10 is a real date/time serial with the same style as the other
# creation_date cells (style index 2).
$ws.Range("G10").Value = 45810.72437689169
$ws.Range("G10").NumberFormat = $ws.Range("G9").NumberFormat
